# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the newer snapshot of data (gh-pages output update).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row number (F column) -> new value
$updates = @{
    2 = 1320
    3 = 1758
    4 = 70
    6 = 6275
    7 = 119
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
